$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking values
# (e.g. "244.07", "1.000") are stored as text, matching the source data
# which uses European-style grouped numbers as plain strings.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '29.907.86'

# Row 3
$ws.Range("D3").Value = '1.895.66'
$ws.Range("E3").Value = '  +0.29%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").Value = '0.7824'
$ws.Range("E5").Value = '  +0.12%  '

# Row 6
$ws.Range("D6").Value = '244.07'
$ws.Range("E6").Value = '  +0.89%  '

# Row 7
$ws.Range("E7").Value = '  +0.06%  '

# Row 8
$ws.Range("D8").Value = '0.3133'
$ws.Range("E8").Value = '  -1.05%  '

# Row 9
$ws.Range("D9").Value = '25.75'
$ws.Range("E9").Value = '  +1.48%  '

# Row 10
$ws.Range("D10").Value = '0.07348'
$ws.Range("E10").Value = '  +4.61%  '

# Row 11
$ws.Range("D11").Value = '0.08092'
$ws.Range("E11").Value = '  +0.63%  '

# Row 12
$ws.Range("D12").Value = '0.7739'

# Row 13
$ws.Range("D13").Value = '5.514'
$ws.Range("E13").Value = '  +4.49%  '

# Row 14
$ws.Range("D14").Value = '1.924.88'
$ws.Range("E14").Value = '  +1.92%  '

# Row 15
$ws.Range("D15").Value = '93.92'
$ws.Range("E15").Value = '  +2.11%  '

# Row 16
$ws.Range("D16").Value = '6.256'
$ws.Range("E16").Value = '  +5.99%  '

# Row 17
$ws.Range("D17").Value = '29.918.16'
$ws.Range("E17").Value = '  +0.49%  '

# Row 18
$ws.Range("D18").Value = '13.99'
$ws.Range("E18").Value = '  +1.12%  '

# Row 19
$ws.Range("D19").Value = '247.34'
$ws.Range("E19").Value = '  +1.83%  '

# Row 20
$ws.Range("D20").Value = '0.000007835'
$ws.Range("E20").Value = '  +1.76%  '

# Row 21
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.157.33'
$ws.Range("E21").Value = '  +0.78%  '

# Row 22
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.05%  '

# Row 23
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").Value = '8.126'
$ws.Range("E23").Value = '  -1.11%  '

# Row 24
$ws.Range("E24").Value = '  +0.03%  '

# Row 25
$ws.Range("D25").Value = '0.1595'
$ws.Range("E25").Value = '  -3.22%  '

# Row 26
$ws.Range("D26").Value = '9.459'
$ws.Range("E26").Value = '  +1.77%  '

# Row 27
$ws.Range("D27").Value = '163.56'
$ws.Range("E27").Value = '  -1.05%  '

# Row 28
$ws.Range("E28").Value = '  +0.36%  '

# Row 29
$ws.Range("D29").Value = '2.026'
$ws.Range("E29").Value = '  -1.06%  '

# Row 30
$ws.Range("D30").Value = '1.436'
$ws.Range("E30").Value = '  +2.93%  '

# Row 31
$ws.Range("D31").Value = '1.544'
$ws.Range("E31").Value = '  +0.68%  '

# Row 32
$ws.Range("D32").Value = '4.492'
$ws.Range("E32").Value = '  +1.96%  '

# Row 33
$ws.Range("E33").Value = '  -0.76%  '

# Row 34
$ws.Range("D34").Value = '4.066'
$ws.Range("E34").Value = '  +0.82%  '

# Row 35
$ws.Range("E35").Value = '  -1.46%  '

# Row 36
$ws.Range("D36").Value = '0.7542'
$ws.Range("E36").Value = '  +2.52%  '

# Row 37
$ws.Range("E37").Value = '  -0.12%  '

# Row 38
$ws.Range("D38").Value = '2.683'
$ws.Range("E38").Value = '  +1.59%  '

# Row 39
$ws.Range("D39").Value = '0.01937'
$ws.Range("E39").Value = '  +1.52%  '

# Row 40
$ws.Range("D40").Value = '2.802'
$ws.Range("E40").Value = '  +1.41%  '

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.4480'
$ws.Range("E41").Value = '  +1.97%  '

# Row 42
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '74.53'
$ws.Range("E42").Value = '  +2.96%  '

# Row 43
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '1.117.06'
$ws.Range("E43").Value = '  +9.64%  '

# Row 44
$ws.Range("D44").Value = '5.965'
$ws.Range("E44").Value = '  +2.76%  '

# Row 45
$ws.Range("D45").Value = '0.8516'
$ws.Range("E45").Value = '  +1.54%  '

# Row 46
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  +0.12%  '

# Row 47
$ws.Range("D47").Value = '1.891'
$ws.Range("E47").Value = '  +1.54%  '

# Row 48
$ws.Range("D48").Value = '102.59'
$ws.Range("E48").Value = '  +0.23%  '

# Row 49
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").Value = '7.529'
$ws.Range("E49").Value = '  +1.79%  '

# Row 50
$ws.Range("B50").Value = 'SynthetixNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D50").Value = '3.053'
$ws.Range("E50").Value = '  +5.30%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '9.765'
$ws.Range("E51").Value = '  -1.17%  '
